$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Fix 1: title said "Assignment II: CUDA Basics" but this is actually
# Assignment III -> insert an "I" so it reads "Assignment III: CUDA Basics".
# ---------------------------------------------------------------------------

# Locate the exact title run so this keeps working even if earlier content
# in the document shifts character offsets around.
$titleFind = $d.Content.Duplicate
$titleFind.Find.Execute("Assignment II: CUDA Basics", $true, $false, $false,
                         $false, $false, $true, 1, $false, "", 0) | Out-Null
$titleStart = $titleFind.Start

# Insert the missing "I" right after "Assignment " (11 characters in).
$insertAt = $titleStart + 11
$insRange = $d.Range($insertAt, $insertAt)
$insRange.InsertBefore("I") | Out-Null

# The engine collapses a touched paragraph's runs down to the minimal set of
# distinct-formatting runs, so after the insert the whole title collapses
# into a single run. Nudge each logical piece (the freshly split "I"/"II: "
# as well as every run that existed before the edit) through a harmless
# Bold-off/Bold-on round trip so the paragraph keeps the same run
# boundaries it originally had, plus the new one introduced by the typed
# "I" - exactly like Word does when you type inside existing text instead
# of doing a find/replace over the whole run.

$p0 = $titleStart + 0
$p1 = $titleStart + 11
$p2 = $titleStart + 12
$p3 = $titleStart + 16
$p4 = $titleStart + 17
$p5 = $titleStart + 18
$p6 = $titleStart + 20
$p7 = $titleStart + 21
$p8 = $titleStart + 27

$rng = $d.Range($p0, $p1)    # "Assignment "
$rng.Bold = 0
$rng.Bold = 1

$rng = $d.Range($p1, $p2)    # "I"            (new)
$rng.Bold = 0
$rng.Bold = 1

$rng = $d.Range($p2, $p3)    # "II: "
$rng.Bold = 0
$rng.Bold = 1

$rng = $d.Range($p3, $p4)    # "C"
$rng.Bold = 0
$rng.Bold = 1

$rng = $d.Range($p4, $p5)    # "U"
$rng.Bold = 0
$rng.Bold = 1

$rng = $d.Range($p5, $p6)    # "DA"
$rng.Bold = 0
$rng.Bold = 1

$rng = $d.Range($p6, $p7)    # " "
$rng.Bold = 0
$rng.Bold = 1

$rng = $d.Range($p7, $p8)    # "Basics"
$rng.Bold = 0
$rng.Bold = 1

# ---------------------------------------------------------------------------
# Fix 2: "<num_elements>" and " (but only with shared memory)" were split
# across two runs with identical formatting; merge that text into one run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("<num_elements> (but only with shared memory)",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "<num_elements> (but only with shared memory)", 2) | Out-Null
